$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Insert a new column before column E ("Training parameters") and move the
#    values that used to live in the last column (L, "Training parameters")
#    into the new column E. This matches the re-ordering seen in the diff
#    where the "Training parameters" column moves from the end of the table
#    to right after "validation samples".
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# After the insert, the former column L ("Training parameters" data) now
# lives in column M. Cut it and paste it into the freshly inserted column E.
$ws.Range("M1:M5").Cut($ws.Range("E1"))
$excel.CutCopyMode = $false

# The Cut above overwrote the styles of E1:E3 with the style that used to be
# on the source cells. Restore the styling that belongs to row 1-3 of this
# table (column D carries the same per-row styling pattern we need here).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the now empty column that used to hold "Training parameters" data
# (it is empty after the cut above, so removing it does not discard any
# content still in use elsewhere).
$ws.Columns("M:M").Delete()

# Give the new column E its own width.
$ws.Range("E:E").ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# 2. Row height tweaks.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").RowHeight = 90
$ws.Rows("5:5").RowHeight = 139.5

# ---------------------------------------------------------------------------
# 3. Add the new model row (row 6).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "lstm128_lr1e-5"
$ws.Range("D6").Value = $ws.Range("D5").Value2
$ws.Range("F6").Value = "100-90.000"
$ws.Range("G6").Value = "90.000-100.000"
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 0.43877539889849798
$ws.Range("J6").Value = 0.86418295996080696
$ws.Range("K6").Value = 0.38935984671115798
$ws.Range("L6").Value = 0.86938856406699005

# Use Copy/Paste for the cells whose source text starts with a literal
# apostrophe so the text-prefix character is not swallowed when read back
# through .Value / .Value2.
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B5").Copy($ws.Range("B6"))
$ws.Range("C5").Copy($ws.Range("C6"))
$ws.Range("E3").Copy($ws.Range("E6"))
$excel.CutCopyMode = $false

# Match the row height used for this new row.
$ws.Rows("6:6").RowHeight = 114.75

# ---------------------------------------------------------------------------
# 4. Update the sheet view: select B6 (this also drops the stale
#    topLeftCell="A3" scroll position from the original view).
# ---------------------------------------------------------------------------
$ws.Range("B6").Select()
